$wb = $excel.ActiveWorkbook

# --- Add the new "Alexandra" worksheet ------------------------------------
# Worksheets.Add() inserts before the active sheet by default, so move the
# new sheet to the very end (after "Example" and "ExampleAllCompounds").
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Alexandra"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch a live reference to the sheet now that it has been relocated.
$ws = $wb.Worksheets.Item("Alexandra")

# --- Populate the calibration standards table ------------------------------
$ws.Range("A1").Value = "Cal.Name"
$ws.Range("B1").Value = "Concentration"

$ws.Range("A2").Value = "stnd300"
$ws.Range("B2").Value = 1/300

$ws.Range("A3").Value = "stnd100"
$ws.Range("B3").Value = 1/100

$ws.Range("A4").Value = "stnd30"
$ws.Range("B4").Value = 1/30

$ws.Range("A5").Value = "stnd10"
$ws.Range("B5").Value = 1/10

$ws.Range("A6").Value = "stnd3"
$ws.Range("B6").Value = 1/3

# --- Make "Alexandra" the active / selected tab ----------------------------
$ws.Activate()
$ws.Range("E9").Select()
